$d = $word.ActiveDocument

# Locate the target word "j'utile" inside the sentence.
$r = $d.Content
$r.Find.Execute("j’utile", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$matchStart = $r.Start
$matchEnd = $r.End

# "j’utile" = "j’uti" (5 chars) + "le" (2 chars)
# We want the final text to read "j’utilise" = "j’uti" + "lise"
$tailRange = $d.Range($matchStart + 5, $matchEnd)

# Temporarily mark the tail with distinct formatting so the replacement
# below doesn't get silently re-merged with its identically-formatted
# neighbour; we'll strip the marker formatting again right after.
$tailRange.Font.Bold = $true
$tailRange.Text = "lise"

$newTailEnd = $matchStart + 5 + 4
$markedRange = $d.Range($matchStart + 5, $newTailEnd)
$markedRange.Font.Bold = $false
